$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the sheet (tab name) to reflect the new "through" date
$ws.Name = "Through 2022-09-05"

# 2. Update the row label for September to reflect the new "through" date
$ws.Range("A10").Value = "September (through 09-05)"

# 3. Update September (row 10) per-year values
$ws.Range("I9").Value = 167

$ws.Range("B10").Value = 3
$ws.Range("C10").Value = 9
$ws.Range("D10").Value = 14
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = 11
$ws.Range("G10").Value = 18
$ws.Range("H10").Value = 22
$ws.Range("I10").Value = 28

# 4. Update Total (row 11) per-year values
$ws.Range("B11").Value = 197
$ws.Range("C11").Value = 390
$ws.Range("D11").Value = 565
$ws.Range("E11").Value = 497
$ws.Range("F11").Value = 360
$ws.Range("G11").Value = 802
$ws.Range("H11").Value = 1092
$ws.Range("I11").Value = 1166
